$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 4, shifting existing rows 4-10 down to 6-12
$ws.Range("A4:B5").Insert()

# Fill the new rows with the new programs
$ws.Range("A4").Value = "RWTH_Aachen_Software_System_Engineering"
$ws.Range("B4").Value = "Yes"
$ws.Range("A5").Value = "RWTH_Aachen_Media_Informatics"
$ws.Range("B5").Value = "Yes"

# Widen column A to fit the new longer program names
$ws.Columns("A").ColumnWidth = 25.5

# Update the active selection like the author's final cursor position
$ws.Range("A6").Select()
